# Edit script for log_pcsmote_x_muestra_wdbc_D50_R25_Pentropia.xlsx
# Updates column F (percentil_dist_25) value for rows 2..171 from
# 0.8767737423380013 to 0.8762693424044758, and refreshes the Z column
# (timestamp) with new per-row timestamp values, matching a rerun of the
# PCSMOTE sample logging process.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newF = 0.8762693424044758

$timestamps = @(
    "2025-10-19T23:54:38.379788", "2025-10-19T23:54:38.379788", "2025-10-19T23:54:38.380785", "2025-10-19T23:54:38.380785", "2025-10-19T23:54:38.380785", "2025-10-19T23:54:38.380785", "2025-10-19T23:54:38.380785", "2025-10-19T23:54:38.380785", "2025-10-19T23:54:38.380785", "2025-10-19T23:54:38.381789",
    "2025-10-19T23:54:38.381789", "2025-10-19T23:54:38.381789", "2025-10-19T23:54:38.381789", "2025-10-19T23:54:38.381789", "2025-10-19T23:54:38.381789", "2025-10-19T23:54:38.381789", "2025-10-19T23:54:38.381789", "2025-10-19T23:54:38.381789", "2025-10-19T23:54:38.381789", "2025-10-19T23:54:38.382787",
    "2025-10-19T23:54:38.382787", "2025-10-19T23:54:38.382787", "2025-10-19T23:54:38.382787", "2025-10-19T23:54:38.382787", "2025-10-19T23:54:38.382787", "2025-10-19T23:54:38.382787", "2025-10-19T23:54:38.382787", "2025-10-19T23:54:38.382787", "2025-10-19T23:54:38.383787", "2025-10-19T23:54:38.383787",
    "2025-10-19T23:54:38.383787", "2025-10-19T23:54:38.383787", "2025-10-19T23:54:38.383787", "2025-10-19T23:54:38.383787", "2025-10-19T23:54:38.383787", "2025-10-19T23:54:38.383787", "2025-10-19T23:54:38.383787", "2025-10-19T23:54:38.383787", "2025-10-19T23:54:38.384786", "2025-10-19T23:54:38.384786",
    "2025-10-19T23:54:38.384786", "2025-10-19T23:54:38.384786", "2025-10-19T23:54:38.384786", "2025-10-19T23:54:38.384786", "2025-10-19T23:54:38.384786", "2025-10-19T23:54:38.384786", "2025-10-19T23:54:38.384786", "2025-10-19T23:54:38.385786", "2025-10-19T23:54:38.385786", "2025-10-19T23:54:38.385786",
    "2025-10-19T23:54:38.386811", "2025-10-19T23:54:38.386811", "2025-10-19T23:54:38.386811", "2025-10-19T23:54:38.387801", "2025-10-19T23:54:38.387801", "2025-10-19T23:54:38.388797", "2025-10-19T23:54:38.388797", "2025-10-19T23:54:38.388797", "2025-10-19T23:54:38.388797", "2025-10-19T23:54:38.388797",
    "2025-10-19T23:54:38.388797", "2025-10-19T23:54:38.389790", "2025-10-19T23:54:38.389790", "2025-10-19T23:54:38.389790", "2025-10-19T23:54:38.389790", "2025-10-19T23:54:38.389790", "2025-10-19T23:54:38.389790", "2025-10-19T23:54:38.389790", "2025-10-19T23:54:38.389790", "2025-10-19T23:54:38.389790",
    "2025-10-19T23:54:38.390786", "2025-10-19T23:54:38.390786", "2025-10-19T23:54:38.390786", "2025-10-19T23:54:38.390786", "2025-10-19T23:54:38.390786", "2025-10-19T23:54:38.390786", "2025-10-19T23:54:38.390786", "2025-10-19T23:54:38.390786", "2025-10-19T23:54:38.390786", "2025-10-19T23:54:38.391787",
    "2025-10-19T23:54:38.391787", "2025-10-19T23:54:38.391787", "2025-10-19T23:54:38.391787", "2025-10-19T23:54:38.391787", "2025-10-19T23:54:38.391787", "2025-10-19T23:54:38.391787", "2025-10-19T23:54:38.391787", "2025-10-19T23:54:38.391787", "2025-10-19T23:54:38.392786", "2025-10-19T23:54:38.392786",
    "2025-10-19T23:54:38.392786", "2025-10-19T23:54:38.392786", "2025-10-19T23:54:38.392786", "2025-10-19T23:54:38.392786", "2025-10-19T23:54:38.392786", "2025-10-19T23:54:38.392786", "2025-10-19T23:54:38.392786", "2025-10-19T23:54:38.392786", "2025-10-19T23:54:38.393786", "2025-10-19T23:54:38.393786",
    "2025-10-19T23:54:38.393786", "2025-10-19T23:54:38.393786", "2025-10-19T23:54:38.393786", "2025-10-19T23:54:38.393786", "2025-10-19T23:54:38.393786", "2025-10-19T23:54:38.393786", "2025-10-19T23:54:38.393786", "2025-10-19T23:54:38.393786", "2025-10-19T23:54:38.394892", "2025-10-19T23:54:38.394892",
    "2025-10-19T23:54:38.395426", "2025-10-19T23:54:38.395426", "2025-10-19T23:54:38.395426", "2025-10-19T23:54:38.395426", "2025-10-19T23:54:38.398419", "2025-10-19T23:54:38.399418", "2025-10-19T23:54:38.399418", "2025-10-19T23:54:38.400420", "2025-10-19T23:54:38.400420", "2025-10-19T23:54:38.400420",
    "2025-10-19T23:54:38.400420", "2025-10-19T23:54:38.401425", "2025-10-19T23:54:38.401425", "2025-10-19T23:54:38.401425", "2025-10-19T23:54:38.401425", "2025-10-19T23:54:38.401425", "2025-10-19T23:54:38.404434", "2025-10-19T23:54:38.404434", "2025-10-19T23:54:38.405424", "2025-10-19T23:54:38.405424",
    "2025-10-19T23:54:38.405424", "2025-10-19T23:54:38.405424", "2025-10-19T23:54:38.406420", "2025-10-19T23:54:38.406420", "2025-10-19T23:54:38.406420", "2025-10-19T23:54:38.406420", "2025-10-19T23:54:38.407423", "2025-10-19T23:54:38.407423", "2025-10-19T23:54:38.407423", "2025-10-19T23:54:38.407423",
    "2025-10-19T23:54:38.408421", "2025-10-19T23:54:38.408421", "2025-10-19T23:54:38.408421", "2025-10-19T23:54:38.409424", "2025-10-19T23:54:38.409424", "2025-10-19T23:54:38.409424", "2025-10-19T23:54:38.409424", "2025-10-19T23:54:38.409424", "2025-10-19T23:54:38.410424", "2025-10-19T23:54:38.410424",
    "2025-10-19T23:54:38.410424", "2025-10-19T23:54:38.411420", "2025-10-19T23:54:38.411420", "2025-10-19T23:54:38.411420", "2025-10-19T23:54:38.411420", "2025-10-19T23:54:38.411420", "2025-10-19T23:54:38.412423", "2025-10-19T23:54:38.412423", "2025-10-19T23:54:38.412423", "2025-10-19T23:54:38.412423",
    "2025-10-19T23:54:38.412423", "2025-10-19T23:54:38.413424", "2025-10-19T23:54:38.413424", "2025-10-19T23:54:38.413424", "2025-10-19T23:54:38.413424", "2025-10-19T23:54:38.413424", "2025-10-19T23:54:38.414422", "2025-10-19T23:54:38.414422", "2025-10-19T23:54:38.414422", "2025-10-19T23:54:38.414422"
)

$startRow = 2
$endRow = 171

for ($r = $startRow; $r -le $endRow; $r++) {
    $idx = $r - $startRow
    $ws.Cells.Item($r, 6).Value = $newF
    $ws.Cells.Item($r, 26).Value = $timestamps[$idx]
}

Write-Host "Updated rows $startRow to $endRow"
